# Apply updates to rows 4-12 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 111936777
$ws.Range("B4").Value = 77650
$ws.Range("Q4").Value = 490056
$ws.Range("R4").Value = 7088709

# Row 5
$ws.Range("A5").Value = 111936768
$ws.Range("B5").Value = 90235
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 3298
$ws.Range("F5").Value = "Trådticka"
$ws.Range("G5").Value = "Climacocystis borealis"
$ws.Range("H5").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q5").Value = 490317
$ws.Range("R5").Value = 7088522

# Row 6
$ws.Range("A6").Value = 111936775
$ws.Range("B6").Value = 89567
$ws.Range("E6").Value = 1204
$ws.Range("F6").Value = "Gränsticka"
$ws.Range("G6").Value = "Phellopilus nigrolimitatus"
$ws.Range("H6").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q6").Value = 490380
$ws.Range("R6").Value = 7088379

# Row 7
$ws.Range("A7").Value = 111936780
$ws.Range("B7").Value = 77650
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 489952
$ws.Range("R7").Value = 7088557

# Row 8
$ws.Range("A8").Value = 111936774
$ws.Range("B8").Value = 56446
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("Q8").Value = 490378
$ws.Range("R8").Value = 7088551
$ws.Range("AC8").Value = "hack"

# Row 9
$ws.Range("A9").Value = 111936781
$ws.Range("B9").Value = 89941
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 4217
$ws.Range("F9").Value = "Blodticka"
$ws.Range("G9").Value = "Meruliopsis taxicola"
$ws.Range("H9").Value = "(Pers.:Fr.) Bondartsev"
$ws.Range("Q9").Value = 490315
$ws.Range("R9").Value = 7088552

# Row 10
$ws.Range("A10").Value = 111936779
$ws.Range("B10").Value = 77650
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 490008
$ws.Range("R10").Value = 7088597
$ws.Range("AC10").ClearContents()

# Row 11
$ws.Range("A11").Value = 111936776
$ws.Range("B11").Value = 77650
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("Q11").Value = 490398
$ws.Range("R11").Value = 7088445

# Row 12
$ws.Range("A12").Value = 111936767
$ws.Range("B12").Value = 90235
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 3298
$ws.Range("F12").Value = "Trådticka"
$ws.Range("G12").Value = "Climacocystis borealis"
$ws.Range("H12").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q12").Value = 490377
$ws.Range("R12").Value = 7088412
